# Updated cryptos list on Mon Aug 14 18:31:02 UTC 2023 with GitHub Actions
#
# Note: several "Price" (column D) values look like plain numbers
# (e.g. 0.9996, 9.016, 0.00000000115). Excel auto-converts such literals to
# numeric cells, which would silently reformat/round them. Prefixing the
# literal with a leading single-quote forces Excel to keep the value as text
# (exactly as it was stored in the original workbook) while the apostrophe
# itself is not stored as part of the cell content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.344.57'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").Value = '1.844.24'
$ws.Range("E3").Value = '  -0.09%  '
$ws.Range("D4").Value = '''0.9996'
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = '''240.30'
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("D6").Value = '''0.6275'
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("E7").Value = '  +0.26%  '
$ws.Range("D8").Value = '''0.07400'
$ws.Range("E8").Value = '  -2.32%  '
$ws.Range("D9").Value = '''0.2890'
$ws.Range("E9").Value = '  -0.84%  '
$ws.Range("D10").Value = '''24.78'
$ws.Range("E10").Value = '  +1.62%  '
$ws.Range("D11").Value = '''0.07721'
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("D12").Value = '1.845.36'
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("D13").Value = '''4.985'
$ws.Range("E13").Value = '  -0.35%  '
$ws.Range("D14").Value = '''0.6777'
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("E15").Value = '  -2.91%  '
$ws.Range("D16").Value = '''82.12'
$ws.Range("E16").Value = '  -1.15%  '
$ws.Range("D17").Value = '''6.263'
$ws.Range("E17").Value = '  +2.43%  '
$ws.Range("D18").Value = '29.399.16'
$ws.Range("E18").Value = '  +0.25%  '
$ws.Range("D19").Value = '''228.94'
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("D20").Value = '''12.30'
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D22").Value = '''7.456'
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("E23").Value = '  +0.29%  '
$ws.Range("D24").Value = '''158.70'
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("D25").Value = '''8.472'
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("D26").Value = '''0.1352'
$ws.Range("E26").Value = '  -2.96%  '
$ws.Range("D27").Value = '''17.43'
$ws.Range("E27").Value = '  -1.05%  '
$ws.Range("D28").Value = '''0.06600'
$ws.Range("E28").Value = '  +16.57%  '
$ws.Range("E29").Value = '  +2.31%  '
$ws.Range("D30").Value = '''1.487'
$ws.Range("E30").Value = '  +0.98%  '
$ws.Range("D31").Value = '''4.071'
$ws.Range("E31").Value = '  -1.02%  '
$ws.Range("D32").Value = '''4.065'
$ws.Range("E32").Value = '  +0.61%  '
$ws.Range("D33").Value = '''1.836'
$ws.Range("E33").Value = '  +0.74%  '
$ws.Range("D34").Value = '''1.139'
$ws.Range("E34").Value = '  -1.22%  '
$ws.Range("D35").Value = '''0.6932'
$ws.Range("E35").Value = '  -1.00%  '
$ws.Range("D36").Value = '''2.576'
$ws.Range("E37").Value = '  +1.89%  '
$ws.Range("D38").Value = '''2.820'
$ws.Range("E38").Value = '  +3.97%  '
$ws.Range("D39").Value = '1.244.94'
$ws.Range("E39").Value = '  +0.26%  '
$ws.Range("D40").Value = '''6.786'
$ws.Range("E40").Value = '  +5.96%  '
$ws.Range("D41").Value = '''0.9361'
$ws.Range("E41").Value = '  +3.74%  '
$ws.Range("D42").Value = '''1.000'
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("D43").Value = '2.024.46'
$ws.Range("E43").Value = '  +0.89%  '
$ws.Range("D44").Value = '''100.61'
$ws.Range("E44").Value = '  -0.98%  '
$ws.Range("D45").Value = '''65.70'
$ws.Range("E45").Value = '  +0.15%  '
$ws.Range("D46").Value = '''1.718'
$ws.Range("E46").Value = '  +2.80%  '
$ws.Range("D47").Value = '''7.040'
$ws.Range("E47").Value = '  -1.18%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''9.016'
$ws.Range("E48").Value = '  +0.35%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '''0.00000000115'
$ws.Range("E49").Value = '  +2.23%  '
$ws.Range("D50").Value = '''0.1150'
$ws.Range("E50").Value = '  -1.11%  '
$ws.Range("D51").Value = '''0.3911'
$ws.Range("E51").Value = '  -0.96%  '
